# Update "想去人数" (interest count) values across the workbook sheets.
# These are data-only refresh values (gh-pages regenerated output),
# mirroring the cell-value updates captured in the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1387
$ws.Range("F5").Value = 214
$ws.Range("F8").Value = 522
$ws.Range("F12").Value = 31132
$ws.Range("F13").Value = 6635
$ws.Range("F17").Value = 96
$ws.Range("F19").Value = 75
$ws.Range("F21").Value = 413
$ws.Range("F25").Value = 364
$ws.Range("F26").Value = 410
$ws.Range("F28").Value = 161
$ws.Range("F30").Value = 715
$ws.Range("F33").Value = 681

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1102
$ws.Range("F7").Value = 4307
$ws.Range("F9").Value = 224
$ws.Range("F19").Value = 4273

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 338

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 338
$ws.Range("F4").Value = 1102
$ws.Range("F7").Value = 1387
$ws.Range("F8").Value = 214
$ws.Range("F11").Value = 523
$ws.Range("F19").Value = 224
$ws.Range("F20").Value = 224
$ws.Range("F27").Value = 96
$ws.Range("F29").Value = 75
$ws.Range("F32").Value = 413
$ws.Range("F36").Value = 364
$ws.Range("F37").Value = 410
$ws.Range("F39").Value = 161
$ws.Range("F41").Value = 715
